$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.918.64'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.25%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.917.13'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.56%  '

$ws.Range("E4").Value = '  -0.15%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '592.59'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.03%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.27'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.36%  '

$ws.Range("E7").Value = '  -0.06%  '

$ws.Range("E8").Value = '  +0.90%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.86'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.30%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.144'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.28%  '

$ws.Range("E11").Value = '  -1.64%  '

$ws.Range("E12").Value = '  +0.81%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '33.56'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.23%  '

$ws.Range("E14").Value = '  -0.02%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.398.98'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.43%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '60.897.66'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.26%  '

$ws.Range("E17").Value = '  -0.71%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.919.02'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.43%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '430.44'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.55%  '

$ws.Range("E20").Value = '  -1.15%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.680'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.69%  '

$ws.Range("E22").Value = '  -0.21%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '81.38'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.81%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.93'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.86%  '

$ws.Range("E25").Value = '  -0.27%  '

$ws.Range("E26").Value = '  +0.69%  '

$ws.Range("E27").Value = '  +0.00%  '

$ws.Range("E28").Value = '  +5.98%  '

$ws.Range("E29").Value = '  -0.15%  '

$ws.Range("E30").Value = '  -0.27%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.04'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.89%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '26.38'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.33%  '

$ws.Range("E33").Value = '  +0.95%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0₃0847'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.34%  '

$ws.Range("E35").Value = '  +0.98%  '

$ws.Range("E36").Value = '  +0.24%  '

$ws.Range("E37").Value = '  +1.94%  '

$ws.Range("E38").Value = '  -1.54%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.122'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.71%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.52'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.45%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.288'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.66%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '39.88'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.30%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '374.95'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.48%  '

$ws.Range("E44").Value = '  -0.60%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.699.47'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.47%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '131.23'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.54%  '

$ws.Range("E48").Value = '  -5.40%  '

$ws.Range("E49").Value = '  +0.15%  '

$ws.Range("E50").Value = '  -3.27%  '

$ws.Range("E51").Value = '  +1.59%  '
